# Examples Analysis and Conversion Utilities
# -------------------------------------------
# On the "Processes" sheet, the "description" column (B) and "type" column
# (E) were swapped - "type" now sits right after the key column (B) and
# "description" moved to the end (E). "fuel" (C) and "product" (D) stay put.
#
# The sheet was also left active/selected on column B (the new "type"
# column), which is what made Excel re-autofit B/C/D's widths and drop E's
# explicit width back to the sheet default.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processes")

# --- Swap the contents of columns B (description) and E (type) ----------
$lastRow = 10
for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $eCell = $ws.Cells.Item($r, 5)

    $bVal = $bCell.Value2
    $eVal = $eCell.Value2

    $bCell.Value2 = $eVal
    $eCell.Value2 = $bVal
}

# --- Resize columns B, C and D; column E returns to the default width ----
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 8.333333333333334

# --- Make "Processes" the active sheet with B1:B10 selected --------------
$ws.Activate()
$ws.Range("B1:B10").Select()
